$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A117").Formula = '="2025-07-29"'
$ws.Range("A117").Value = $ws.Range("A117").Value
$ws.Range("B117").Value = "ABB"
$ws.Range("C117").Value = "Independiente Petrolero"
$ws.Range("D117").Value = 1
$ws.Range("E117").Value = 2
$ws.Range("F117").Value = 1378247
$ws.Range("G117").Value = 3
$ws.Range("H117").Value = 1
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 2
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 1
$ws.Range("M117").Value = 0
$ws.Range("N117").Value = 2
$ws.Range("O117").Value = 1
$ws.Range("P117").Value = 0
$ws.Range("Q117").Value = 50
$ws.Range("R117").Value = 50
$ws.Range("S117").Value = "V"

$ws.Range("A118").Formula = '="2025-07-30"'
$ws.Range("A118").Value = $ws.Range("A118").Value
$ws.Range("B118").Value = "Nacional Potosí"
$ws.Range("C118").Value = "Oriente Petrolero"
$ws.Range("D118").Value = 4
$ws.Range("E118").Value = 1
$ws.Range("F118").Value = 1370697
$ws.Range("G118").Value = 3
$ws.Range("H118").Value = 3
$ws.Range("I118").Value = 4
$ws.Range("J118").Value = 4
$ws.Range("K118").Value = 1
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 2
$ws.Range("N118").Value = 0
$ws.Range("O118").Value = 2
$ws.Range("P118").Value = 1
$ws.Range("Q118").Value = 42
$ws.Range("R118").Value = 58
$ws.Range("S118").Value = "L"
